$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the two header cells (L1, M1) to their new, more descriptive names.
# (Set M1 before L1 so the shared-string table gets the two new entries in
# the same order as the reference workbook.)
$ws.Range("M1").Value = "date_of_archiving"
$ws.Range("L1").Value = "date_of_data_provision"

# Excel auto-resizes these two columns once the headers change (no longer
# marked as "best fit" — the widths below match (as closely as the engine's
# ColumnWidth rounding allows) the saved widths after the header text change).
$ws.Columns.Item(12).ColumnWidth = 21.67
$ws.Columns.Item(13).ColumnWidth = 17

# Update the view state: scroll so column G is the left-most visible column,
# and move the active selection to L3.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 7
$ws.Range("L3").Select()
